$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "PERIODE_BULANAN" / "VERIFIKASI" values for the regression date update
$ws.Range("Q2").Value = "202405"
$ws.Range("T2").Value = "15/04/2024"

# Update the saved view state (scroll position / selection)
$excel.ActiveWindow.ScrollColumn = 10
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("R2").Select()
